$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are scraped as literal text; force Text format before
# assigning any value that Excel would otherwise auto-convert to a number
# (this avoids turning e.g. "245.04" or "0.05910" into numeric 245.04 / 0.0591,
# which would silently drop the trailing zero / string formatting).

$ws.Range("D2").Value = "30.120.53"
$ws.Range("E2").Value = "  -4.58%  "
$ws.Range("D3").Value = "1.911.77"
$ws.Range("E3").Value = "  -4.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.04"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6973"
$ws.Range("E6").Value = "  -13.83%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3207"
$ws.Range("E8").Value = "  -6.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.60"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06806"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7838"
$ws.Range("E11").Value = "  -7.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07934"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "1.914.75"
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.349"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.34"
$ws.Range("E15").Value = "  -9.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "259.08"
$ws.Range("E16").Value = "  -6.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.26"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "30.133.20"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.771"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007797"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "2.172.62"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.764"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.495"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.11"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.63"
$ws.Range("E27").Value = "  -6.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1297"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.203"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.545"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.371"
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.151"
$ws.Range("E33").Value = "  -4.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05005"
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7358"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01906"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.792"
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.97"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.473"
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4380"
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.997"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8284"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.42"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.701"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.151"
$ws.Range("E48").Value = "  -5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.75"
$ws.Range("E49").Value = "  -2.27%  "

# Row 50/51 swap (NEARProtocol <-> Cronos)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05910"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.461"
$ws.Range("E51").Value = "  +1.30%  "
